$d = $word.ActiveDocument

# The new bullet belongs under the "Principal Cloud and Platform Engineer"
# role, as the first bullet point (immediately before the existing
# "Managed the Base AMI..." bullet). Anchor the insertion on that existing
# bullet's text so the new paragraph lands in the right spot and -- by
# inheriting that paragraph's formatting when we split it -- automatically
# picks up the same "Compact" style and numbered-list (numId 1002)
# formatting used by the rest of the bullets in this job entry.
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "Managed the Base",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor text 'Managed the Base' in the document."
}

$insertPoint = $anchor.Duplicate
$insertPoint.Collapse(1)

# Insert a paragraph break right before the anchor bullet; the new empty
# paragraph created takes on the same paragraph formatting (Compact style,
# numId 1002 bullet numbering) as the paragraph it was split from.
$insertPoint.InsertParagraphBefore()

# Fill in the new (now-empty) paragraph with the new bullet's text.
$insertPoint.InsertBefore("As every school in America transitioned to online learning during the COVID-19 lockdowns, I was the technical/development lead on the team who supported all SRE and product engineering teams, working on core platforms and services.")

Write-Output "Inserted new bullet before 'Managed the Base...' paragraph."
